$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3432399.5
$ws.Range("I9").Value = 8000
$ws.Range("K9").Value = 8000
$ws.Range("M9").Value = -7831
$ws.Range("H11").Value = 4215.8887
$ws.Range("I11").Value = 4215.8887
$ws.Range("K11").Value = 4215.8887
$ws.Range("M11").Value = -4075.8887
$ws.Range("H12").Value = 403.25
$ws.Range("I12").Value = 337.66666
$ws.Range("K12").Value = 337.66666
$ws.Range("M12").Value = -167.66666
$ws.Range("H15").Value = 2709700.8
$ws.Range("I15").Value = 2709700.8
$ws.Range("K15").Value = 8129102.399999999
$ws.Range("M15").Value = -8128933.399999999
$ws.Range("H40").Value = 2916.25
$ws.Range("I40").Value = 2545
$ws.Range("K40").Value = 2545
$ws.Range("M40").Value = -2370
$ws.Range("H70").Value = 5952.4375
$ws.Range("I70").Value = 5166.8335
$ws.Range("J70").Value = 6423.8
$ws.Range("K70").Value = 15500.5005
$ws.Range("L70").Value = 19271.4
$ws.Range("M70").Value = -15230.5005
$ws.Range("N70").Value = -19811.4
$ws.Range("H73").Value = 5952.4375
$ws.Range("I73").Value = 5166.8335
$ws.Range("J73").Value = 6423.8
$ws.Range("K73").Value = 15500.5005
$ws.Range("L73").Value = 19271.4
$ws.Range("M73").Value = -14564.5005
$ws.Range("N73").Value = -21143.4
$ws.Range("H86").Value = 1850
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -3446
$ws.Range("H89").Value = 1850
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -17232
$ws.Range("H132").Value = 1236.2916
$ws.Range("I132").Value = 1098.75
$ws.Range("K132").Value = 3296.25
$ws.Range("M132").Value = -766.25
$ws.Range("H138").Value = 3647.257
$ws.Range("I138").Value = 2773.3333
$ws.Range("K138").Value = 8319.999899999999
$ws.Range("M138").Value = -3179.999899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21084.074
$ws.Range("I32").Value = 21084.074
$ws.Range("K32").Value = 21084.074
$ws.Range("M32").Value = -20797.074
$ws.Range("H132").Value = 2554.58
$ws.Range("I132").Value = 1999.3422
$ws.Range("K132").Value = 5998.0266
$ws.Range("M132").Value = -3468.0266

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4251.5713
$ws.Range("J20").Value = 3808.125
$ws.Range("L20").Value = 3808.125
$ws.Range("N20").Value = -4302.125
$ws.Range("H22").Value = 1321.625
$ws.Range("I22").Value = 1474.7142
$ws.Range("K22").Value = 1474.7142
$ws.Range("M22").Value = -1301.7142
$ws.Range("H86").Value = 338212.66
$ws.Range("I86").Value = 4295.8335
$ws.Range("J86").Value = 672129.5
$ws.Range("K86").Value = 4295.8335
$ws.Range("L86").Value = 672129.5
$ws.Range("M86").Value = -3172.8335
$ws.Range("N86").Value = -674375.5
$ws.Range("H89").Value = 338212.66
$ws.Range("I89").Value = 4295.8335
$ws.Range("J89").Value = 672129.5
$ws.Range("K89").Value = 21479.1675
$ws.Range("L89").Value = 3360647.5
$ws.Range("M89").Value = -15863.1675
$ws.Range("N89").Value = -3371879.5
$ws.Range("H94").Value = 2574.75
$ws.Range("I94").Value = 2574.75
$ws.Range("K94").Value = 2574.75
$ws.Range("M94").Value = -2123.75
$ws.Range("H105").Value = 100026456
$ws.Range("I105").Value = 100026456
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 100026456
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -100024709
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 3685.2
$ws.Range("I107").Value = 3540.1667
$ws.Range("J107").Value = 4265.3335
$ws.Range("K107").Value = 3540.1667
$ws.Range("L107").Value = 4265.3335
$ws.Range("M107").Value = -1620.1667
$ws.Range("N107").Value = -8105.3335

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5204.8887
$ws.Range("I16").Value = 4209.4
$ws.Range("K16").Value = 4209.4
$ws.Range("M16").Value = -3922.4
$ws.Range("H22").Value = 633.8421
$ws.Range("I22").Value = 587.7143
$ws.Range("K22").Value = 587.7143
$ws.Range("M22").Value = -237.7143
$ws.Range("H31").Value = 25002548
$ws.Range("I31").Value = 47619784
$ws.Range("J31").Value = 4547.263
$ws.Range("K31").Value = 47619784
$ws.Range("L31").Value = 4547.263
$ws.Range("M31").Value = -47619489
$ws.Range("N31").Value = -5137.263
$ws.Range("H34").Value = 25002548
$ws.Range("I34").Value = 47619784
$ws.Range("J34").Value = 4547.263
$ws.Range("K34").Value = 47619784
$ws.Range("L34").Value = 4547.263
$ws.Range("M34").Value = -47619582
$ws.Range("N34").Value = -4951.263
$ws.Range("H113").Value = 5204.8887
$ws.Range("I113").Value = 4209.4
$ws.Range("K113").Value = 4209.4
$ws.Range("M113").Value = -2039.4
$ws.Range("H132").Value = 5580.85
$ws.Range("I132").Value = 4113.273
$ws.Range("J132").Value = 7374.5557
$ws.Range("K132").Value = 12339.819
$ws.Range("L132").Value = 22123.6671
$ws.Range("M132").Value = -9809.819
$ws.Range("N132").Value = -27183.6671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50973260
$ws.Range("I4").Value = 80900136
$ws.Range("J4").Value = 7745552
$ws.Range("K4").Value = 242700408
$ws.Range("L4").Value = 23236656
$ws.Range("M4").Value = -242700296
$ws.Range("N4").Value = -23236880
$ws.Range("H33").Value = 121.35714
$ws.Range("I33").Value = 105
$ws.Range("K33").Value = 630
$ws.Range("M33").Value = -347
$ws.Range("H49").Value = 900
$ws.Range("J49").Value = 500
$ws.Range("L49").Value = 1500
$ws.Range("N49").Value = -1812

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 72662.07000000001
$ws.Range("I2").Value = 1428.9166
$ws.Range("K2").Value = 1428.9166
$ws.Range("M2").Value = -1315.9166
$ws.Range("H70").Value = 17243
$ws.Range("J70").Value = 13435.5
$ws.Range("L70").Value = 13435.5
$ws.Range("N70").Value = -13975.5
$ws.Range("H73").Value = 17243
$ws.Range("J73").Value = 13435.5
$ws.Range("L73").Value = 13435.5
$ws.Range("N73").Value = -15307.5
$ws.Range("H102").Value = 1244.6923
$ws.Range("I102").Value = 825.5
$ws.Range("K102").Value = 825.5
$ws.Range("M102").Value = 796.5
$ws.Range("H132").Value = 5912.143
$ws.Range("I132").Value = 3439.7778
$ws.Range("K132").Value = 10319.3334
$ws.Range("M132").Value = -7789.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1850.3889
$ws.Range("I16").Value = 1387.1333
$ws.Range("K16").Value = 1387.1333
$ws.Range("M16").Value = -1217.1333
$ws.Range("H46").Value = 4196.9756
$ws.Range("I46").Value = 1781.1875
$ws.Range("K46").Value = 1781.1875
$ws.Range("M46").Value = -1593.1875
$ws.Range("H61").Value = 4024.0356
$ws.Range("I61").Value = 4099.185
$ws.Range("K61").Value = 4099.185
$ws.Range("M61").Value = -3897.185
$ws.Range("H113").Value = 4024.0356
$ws.Range("I113").Value = 4099.185
$ws.Range("K113").Value = 4099.185
$ws.Range("M113").Value = -1929.185
$ws.Range("H136").Value = 4399.619
$ws.Range("I136").Value = 1714.7778
$ws.Range("J136").Value = 6413.25
$ws.Range("K136").Value = 5144.3334
$ws.Range("L136").Value = 19239.75
$ws.Range("M136").Value = -2594.3334
$ws.Range("N136").Value = -24339.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1642.8334
$ws.Range("I107").Value = 1445.3529
$ws.Range("K107").Value = 4336.0587
$ws.Range("M107").Value = -2416.0587
$ws.Range("H132").Value = 7765.375
$ws.Range("I132").Value = 5772.5835
$ws.Range("K132").Value = 17317.7505
$ws.Range("M132").Value = -14787.7505
